# Update simulation results for the 380 kV case (res_line/pl_mw sheet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1803039715586294
$ws.Range("C2").Value = 0.04861204301444388
$ws.Range("D2").Value = 0.02683154439476709
$ws.Range("F2").Value = 0.4712653494825716
$ws.Range("G2").Value = 0.3148194623067724
$ws.Range("H2").Value = 0.4904941423470817
$ws.Range("I2").Value = 0.3830886459092078
$ws.Range("K2").Value = 0.1899372257878866
$ws.Range("M2").Value = 0.8901613683944163
$ws.Range("N2").Value = 1.188589985725999
$ws.Range("O2").Value = 1.540774767514691

# Row 3
$ws.Range("B3").Value = 0.1572978465636368
$ws.Range("C3").Value = 0.04643923291091312
$ws.Range("D3").Value = 0.02382316355171099
$ws.Range("F3").Value = 0.4700020403240401
$ws.Range("G3").Value = 0.3149087262283672
$ws.Range("H3").Value = 0.4934004564138093
$ws.Range("I3").Value = 0.3862387656581241
$ws.Range("K3").Value = 0.166854195485115
$ws.Range("M3").Value = 0.792407016919725
$ws.Range("N3").Value = 1.202878906655356
$ws.Range("O3").Value = 1.546843590894994

# Row 4
$ws.Range("B4").Value = 0.143119466498149
$ws.Range("C4").Value = 0.04509321617330642
$ws.Range("D4").Value = 0.02196411058593384
$ws.Range("F4").Value = 0.469510832096276
$ws.Range("G4").Value = 0.3151799816537491
$ws.Range("H4").Value = 0.4953793536497102
$ws.Range("I4").Value = 0.3883732312516415
$ws.Range("K4").Value = 0.1526028471889163
$ws.Range("M4").Value = 0.7328981293851626
$ws.Range("N4").Value = 1.212082274134679
$ws.Range("O4").Value = 1.551428629058989

# Row 5
$ws.Range("B5").Value = 0.1373289744420134
$ws.Range("C5").Value = 0.04454176711711
$ws.Range("D5").Value = 0.0212035890839033
$ws.Range("F5").Value = 0.4693822208720633
$ws.Range("G5").Value = 0.3153448979572673
$ws.Range("H5").Value = 0.4962347002729857
$ws.Range("I5").Value = 0.3892933893658324
$ws.Range("K5").Value = 0.1467761439531898
$ws.Range("M5").Value = 0.7087709747411424
$ws.Range("N5").Value = 1.215940798823877
$ws.Range("O5").Value = 1.553512992369534

# Row 6
$ws.Range("B6").Value = 0.1363667167054814
$ws.Range("C6").Value = 0.04445002374188789
$ws.Range("D6").Value = 0.0210771287904592
$ws.Range("F6").Value = 0.4693651879663534
$ws.Range("G6").Value = 0.3153755649992647
$ws.Range("H6").Value = 0.4963796864822285
$ws.Range("I6").Value = 0.3894492212405076
$ws.Range("K6").Value = 0.1458074821937743
$ws.Range("M6").Value = 0.7047719548222489
$ws.Range("N6").Value = 1.216588030779114
$ws.Range("O6").Value = 1.553872140030052

# Row 7
$ws.Range("B7").Value = 0.1430414245760545
$ws.Range("C7").Value = 0.04508579095237053
$ws.Range("D7").Value = 0.02195386577597702
$ws.Range("F7").Value = 0.4695088078201479
$ws.Range("G7").Value = 0.3151819856674365
$ws.Range("H7").Value = 0.4953906909808055
$ws.Range("I7").Value = 0.3883854369992115
$ws.Range("K7").Value = 0.1525243431346865
$ws.Range("M7").Value = 0.7325722507725629
$ws.Range("N7").Value = 1.212133873991371
$ws.Range("O7").Value = 1.551455865277063

# Row 8
$ws.Range("B8").Value = 0.1723826672283337
$ws.Range("C8").Value = 0.04786536187861401
$ws.Range("D8").Value = 0.02579675193602071
$ws.Range("F8").Value = 0.4707707340662282
$ws.Range("G8").Value = 0.3148052859046757
$ws.Range("H8").Value = 0.4914559162597811
$ws.Range("I8").Value = 0.3841332295565962
$ws.Range("K8").Value = 0.181994747786888
$ws.Range("M8").Value = 0.8563457637780942
$ws.Range("N8").Value = 1.193427562215492
$ws.Range("O8").Value = 1.542689049208533

# Row 9
$ws.Range("B9").Value = 0.2294836828769178
$ws.Range("C9").Value = 0.05321938129580417
$ws.Range("D9").Value = 0.0332364369936613
$ws.Range("F9").Value = 0.4755021080044131
$ws.Range("G9").Value = 0.3157867250935666
$ws.Range("H9").Value = 0.485280553293741
$ws.Range("I9").Value = 0.3773846078059186
$ws.Range("K9").Value = 0.2391447538212645
$ws.Range("M9").Value = 1.103406273525096
$ws.Range("N9").Value = 1.160156340792336
$ws.Range("O9").Value = 1.532313407733625

# Row 10
$ws.Range("B10").Value = 0.2711458738884289
$ws.Range("C10").Value = 0.05709129498663401
$ws.Range("D10").Value = 0.03864174833454115
$ws.Range("F10").Value = 0.4803547481183941
$ws.Range("G10").Value = 0.3175609297515152
$ws.Range("H10").Value = 0.4816804984546508
$ws.Range("I10").Value = 0.373396952979359
$ws.Range("K10").Value = 0.2807183062101331
$ws.Range("M10").Value = 1.287968666403913
$ws.Range("N10").Value = 1.137791935964764
$ws.Range("O10").Value = 1.528850765033326

# Row 11
$ws.Range("B11").Value = 0.290031366616347
$ws.Range("C11").Value = 0.0588387721519581
$ws.Range("D11").Value = 0.0410872039012844
$ws.Range("F11").Value = 0.4828614739843289
$ws.Range("G11").Value = 0.3185977572337748
$ws.Range("H11").Value = 0.4802457362471202
$ws.Range("I11").Value = 0.3717939108569084
$ws.Range("K11").Value = 0.2995363840598202
$ws.Range("M11").Value = 1.372682297476132
$ws.Range("N11").Value = 1.12806945843649
$ws.Range("O11").Value = 1.528180157315589

# Row 12
$ws.Range("B12").Value = 0.297172711300874
$ws.Range("C12").Value = 0.05949844594790932
$ws.Range("D12").Value = 0.04201125543572459
$ws.Range("F12").Value = 0.4838537306720454
$ws.Range("G12").Value = 0.3190234819181583
$ws.Range("H12").Value = 0.4797315707602223
$ws.Range("I12").Value = 0.3712172417056756
$ws.Range("K12").Value = 0.3066483033856571
$ws.Range("M12").Value = 1.404876915242482
$ws.Range("N12").Value = 1.124452730294333
$ws.Range("O12").Value = 1.528056364599919

# Row 13
$ws.Range("B13").Value = 0.2956351568346918
$ws.Range("C13").Value = 0.05935646591478871
$ws.Range("D13").Value = 0.04181233377290994
$ws.Range("F13").Value = 0.483638117632502
$ws.Range("G13").Value = 0.3189303212899262
$ws.Range("H13").Value = 0.4798410096451562
$ws.Range("I13").Value = 0.371340086581462
$ws.Range("K13").Value = 0.305117259525332
$ws.Range("M13").Value = 1.397937985645257
$ws.Range("N13").Value = 1.125228766839308
$ws.Range("O13").Value = 1.528077235808865

# Row 14
$ws.Range("B14").Value = 0.2906190966547513
$ws.Range("C14").Value = 0.05889308551269323
$ws.Range("D14").Value = 0.04116326636791712
$ws.Range("F14").Value = 0.4829422455650416
$ws.Range("G14").Value = 0.3186321180999983
$ws.Range("H14").Value = 0.4802028515593122
$ws.Range("I14").Value = 0.3717458591320799
$ws.Range("K14").Value = 0.3001217713595281
$ws.Range("M14").Value = 1.375328615655434
$ws.Range("N14").Value = 1.127770605605974
$ws.Range("O14").Value = 1.528167363936689

# Row 15
$ws.Range("B15").Value = 0.2875452721317231
$ws.Range("C15").Value = 0.05860898178228524
$ws.Range("D15").Value = 0.04076543310465297
$ws.Range("F15").Value = 0.4825216050551191
$ws.Range("G15").Value = 0.3184537726529442
$ws.Range("H15").Value = 0.4804282851442423
$ws.Range("I15").Value = 0.371998362315928
$ws.Range("K15").Value = 0.2970600402807406
$ws.Range("M15").Value = 1.361494957307272
$ws.Range("N15").Value = 1.129336017802938
$ws.Range("O15").Value = 1.52823952191487

# Row 16
$ws.Range("B16").Value = 0.2699102687499249
$ws.Range("C16").Value = 0.05697680868217958
$ws.Range("D16").Value = 0.03848165670878245
$ws.Range("F16").Value = 0.4801969477334396
$ws.Range("G16").Value = 0.3174977990307042
$ws.Range("H16").Value = 0.4817783440862584
$ws.Range("I16").Value = 0.3735059633776849
$ws.Range("K16").Value = 0.2794865615437061
$ws.Range("M16").Value = 1.282448206668079
$ws.Range("N16").Value = 1.138436408913677
$ws.Range("O16").Value = 1.528912799512995

# Row 17
$ws.Range("B17").Value = 0.2590742393745415
$ws.Range("C17").Value = 0.0559719249505406
$ws.Range("D17").Value = 0.03707715006945023
$ws.Range("F17").Value = 0.4788474760391424
$ws.Range("G17").Value = 0.3169702267142753
$ws.Range("H17").Value = 0.4826585115676494
$ws.Range("I17").Value = 0.374484882917276
$ws.Range("K17").Value = 0.2686813358698998
$ws.Range("M17").Value = 1.234153813242301
$ws.Range("N17").Value = 1.144134824321421
$ws.Range("O17").Value = 1.5295575685915

# Row 18
$ws.Range("B18").Value = 0.2528353839391002
$ws.Range("C18").Value = 0.05539264009779288
$ws.Range("D18").Value = 0.03626805258093668
$ws.Range("F18").Value = 0.4780994620358499
$ws.Range("G18").Value = 0.3166884007145256
$ws.Range("H18").Value = 0.4831838623835978
$ws.Range("I18").Value = 0.3750677863638181
$ws.Range("K18").Value = 0.2624576478785059
$ws.Range("M18").Value = 1.206446828938908
$ws.Range("N18").Value = 1.14745486644809
$ws.Range("O18").Value = 1.530013560132346

# Row 19
$ws.Range("B19").Value = 0.2507219572191559
$ws.Range("C19").Value = 0.05519628243067132
$ws.Range("D19").Value = 0.03599389094453898
$ws.Range("F19").Value = 0.4778510354299499
$ws.Range("G19").Value = 0.3165966904773754
$ws.Range("H19").Value = 0.4833650188648306
$ws.Range("I19").Value = 0.3752685565665423
$ws.Range("K19").Value = 0.2603489186491998
$ws.Range("M19").Value = 1.19707765396052
$ws.Range("N19").Value = 1.148586266531737
$ws.Range("O19").Value = 1.53018257114293

# Row 20
$ws.Range("B20").Value = 0.2602284049946491
$ws.Range("C20").Value = 0.05607903177551066
$ws.Range("D20").Value = 0.0372267932563517
$ws.Range("F20").Value = 0.4789882146008679
$ws.Range("G20").Value = 0.3170241496649737
$ws.Range("H20").Value = 0.4825628395098605
$ws.Range("I20").Value = 0.374378620132525
$ws.Range("K20").Value = 0.2698324863079904
$ws.Range("M20").Value = 1.239287470721706
$ws.Range("N20").Value = 1.143523822805975
$ws.Range("O20").Value = 1.529480119874563

# Row 21
$ws.Range("B21").Value = 0.292092715979237
$ws.Range("C21").Value = 0.0590292479199519
$ws.Range("D21").Value = 0.04135396763558674
$ws.Range("F21").Value = 0.4831454727975952
$ws.Range("G21").Value = 0.3187188087429647
$ws.Range("H21").Value = 0.4800957790072573
$ws.Range("I21").Value = 0.3716258494591393
$ws.Range("K21").Value = 0.301589453951749
$ws.Range("M21").Value = 1.381966348652554
$ws.Range("N21").Value = 1.127022242206309
$ws.Range("O21").Value = 1.528137358222722

# Row 22
$ws.Range("B22").Value = 0.3128582753230376
$ws.Range("C22").Value = 0.0609453606960102
$ws.Range("D22").Value = 0.0440396928296849
$ws.Range("F22").Value = 0.4861131924985784
$ws.Range("G22").Value = 0.3200193325426568
$ws.Range("H22").Value = 0.478653293579157
$ws.Range("I22").Value = 0.3700037783771393
$ws.Range("K22").Value = 0.3222621440703222
$ws.Range("M22").Value = 1.475891181505602
$ws.Range("N22").Value = 1.116616199574963
$ws.Range("O22").Value = 1.528018446679596

# Row 23
$ws.Range("B23").Value = 0.3017809478189974
$ws.Range("C23").Value = 0.05992381614578335
$ws.Range("D23").Value = 0.04260735260389481
$ws.Range("F23").Value = 0.484506329777922
$ws.Range("G23").Value = 0.319307540489099
$ws.Range("H23").Value = 0.4794076415203889
$ws.Range("I23").Value = 0.3708532998893972
$ws.Range("K23").Value = 0.3112364613189698
$ws.Range("M23").Value = 1.425697592879885
$ws.Range("N23").Value = 1.122135422350219
$ws.Range("O23").Value = 1.528012468994802

# Row 24
$ws.Range("B24").Value = 0.2597066349505894
$ws.Range("C24").Value = 0.0560306136466977
$ws.Range("D24").Value = 0.03715914463799663
$ws.Range("F24").Value = 0.4789245000576017
$ws.Range("G24").Value = 0.3169997041911969
$ws.Range("H24").Value = 0.4826060326333916
$ws.Range("I24").Value = 0.3744265989031028
$ws.Range("K24").Value = 0.2693120873449288
$ws.Range("M24").Value = 1.236966363700688
$ws.Range("N24").Value = 1.143799919711407
$ws.Range("O24").Value = 1.52951486873836

# Row 25
$ws.Range("B25").Value = 0.2140855643677071
$ws.Range("C25").Value = 0.05178161236108281
$ws.Range("D25").Value = 0.03123430781828063
$ws.Range("F25").Value = 0.4739804703693977
$ws.Range("G25").Value = 0.3153366021099373
$ws.Range("H25").Value = 0.4867864449848298
$ws.Range("I25").Value = 0.3790399427761209
$ws.Range("K25").Value = 0.2237552618046124
$ws.Range("M25").Value = 1.036064238498369
$ws.Range("N25").Value = 1.188791839192783
$ws.Range("O25").Value = 1.53439006507476
